$d = $word.ActiveDocument

# 1) Update the "Curso (semestre ideal)" line: drop the "EQD (3), " part.
$d.Content.Find.Execute("Curso (semestre ideal): EQD (3), EQN (3)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Curso (semestre ideal): EQN (3)", 2) | Out-Null

# 2) Append a new "Requisitos" section (Heading2) followed by a bullet item
#    listing the weak prerequisite, mirroring the existing document structure.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null

$heading = $d.Paragraphs.Item($d.Paragraphs.Count)
$heading.Style = "Heading2"
$heading.Range.Text = "Requisitos"

$lastPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r2 = $lastPara2.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter() | Out-Null

$bullet = $d.Paragraphs.Item($d.Paragraphs.Count)
$bullet.Style = "ListBullet"
$bullet.Range.Text = "LOQ4073 -  Química Geral II  (Requisito fraco)"
$bullet.Range.InsertAfter([char]11)
